$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = 28

$ws.Range("E42").Select()
